$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '30.171.93'
$dCell.ClearFormats()
$ws.Range("E2").Value = '  -0.97%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.836.76'
$dCell.ClearFormats()
$ws.Range("E3").Value = '  -1.53%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = '1.000'
$dCell.ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '233.42'
$dCell.ClearFormats()
$ws.Range("E5").Value = '  -0.76%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '1.000'
$dCell.ClearFormats()
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = '0.4670'
$dCell.ClearFormats()
$ws.Range("E7").Value = '  -3.28%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '0.2704'
$dCell.ClearFormats()
$ws.Range("E8").Value = '  -3.56%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.06271'
$dCell.ClearFormats()
$ws.Range("E9").Value = '  -3.64%  '

$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '1.854.48'
$dCell.ClearFormats()
$ws.Range("E10").Value = '  -2.08%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '0.07395'
$dCell.ClearFormats()
$ws.Range("E11").Value = '  -0.48%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '16.03'
$dCell.ClearFormats()
$ws.Range("E12").Value = '  -1.87%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '4.923'
$dCell.ClearFormats()
$ws.Range("E13").Value = '  -2.80%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '83.64'
$dCell.ClearFormats()
$ws.Range("E14").Value = '  -4.04%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '0.6173'
$dCell.ClearFormats()
$ws.Range("E15").Value = '  -4.48%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '30.083.94'
$dCell.ClearFormats()
$ws.Range("E16").Value = '  -1.17%  '

$ws.Range("B17").Value = 'Dai'
$ws.Range("C17").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '1.000'
$dCell.ClearFormats()
$ws.Range("E17").Value = '  -0.02%  '

$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '228.44'
$dCell.ClearFormats()
$ws.Range("E18").Value = '  -2.41%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '0.000007274'
$dCell.ClearFormats()
$ws.Range("E19").Value = '  -3.48%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '12.35'
$dCell.ClearFormats()
$ws.Range("E20").Value = '  -4.79%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '2.078.39'
$dCell.ClearFormats()
$ws.Range("E21").Value = '  -1.68%  '

$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9991'
$dCell.ClearFormats()
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '4.855'
$dCell.ClearFormats()
$ws.Range("E23").Value = '  -5.75%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '5.834'
$dCell.ClearFormats()
$ws.Range("E24").Value = '  -4.17%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '9.186'
$dCell.ClearFormats()
$ws.Range("E25").Value = '  -1.54%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '164.95'
$dCell.ClearFormats()
$ws.Range("E26").Value = '  -1.34%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '17.71'
$dCell.ClearFormats()
$ws.Range("E27").Value = '  -3.49%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '1.874'
$dCell.ClearFormats()
$ws.Range("E28").Value = '  -2.40%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '0.1031'
$dCell.ClearFormats()
$ws.Range("E29").Value = '  +0.41%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '1.370'
$dCell.ClearFormats()
$ws.Range("E30").Value = '  -0.26%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '4.073'
$dCell.ClearFormats()
$ws.Range("E31").Value = '  -4.53%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '3.788'
$dCell.ClearFormats()
$ws.Range("E32").Value = '  -5.20%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '0.04783'
$dCell.ClearFormats()
$ws.Range("E33").Value = '  -3.92%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '1.134'
$dCell.ClearFormats()
$ws.Range("E34").Value = '  -3.71%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '0.7091'
$dCell.ClearFormats()
$ws.Range("E35").Value = '  -2.97%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '2.707'
$dCell.ClearFormats()
$ws.Range("E36").Value = '  -0.25%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '0.01867'
$dCell.ClearFormats()
$ws.Range("E37").Value = '  -2.39%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '2.646'
$dCell.ClearFormats()
$ws.Range("E38").Value = '  +0.59%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '0.8933'
$dCell.ClearFormats()
$ws.Range("E39").Value = '  -2.15%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '1.930'
$dCell.ClearFormats()
$ws.Range("E40").Value = '  -5.60%  '

$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '104.23'
$dCell.ClearFormats()
$ws.Range("E41").Value = '  -2.00%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '1.001'
$dCell.ClearFormats()
$ws.Range("E42").Value = '  +0.56%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '5.519'
$dCell.ClearFormats()
$ws.Range("E43").Value = '  -0.89%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '0.4001'
$dCell.ClearFormats()
$ws.Range("E44").Value = '  -4.72%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '6.946'
$dCell.ClearFormats()
$ws.Range("E45").Value = '  -3.85%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '0.1188'
$dCell.ClearFormats()
$ws.Range("E46").Value = '  -3.21%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '59.71'
$dCell.ClearFormats()
$ws.Range("E47").Value = '  -3.90%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '8.516'
$dCell.ClearFormats()
$ws.Range("E48").Value = '  -3.83%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '32.49'
$dCell.ClearFormats()
$ws.Range("E49").Value = '  -3.25%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '0.05507'
$dCell.ClearFormats()
$ws.Range("E50").Value = '  -2.49%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '1.357'
$dCell.ClearFormats()
$ws.Range("E51").Value = '  -5.83%  '
